$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.393.44"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.17%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.946.62"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.74%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "488.64"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +10.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.19"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.46%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.732"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.59%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.166"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +11.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000353"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +14.82%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.17"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.25%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.589.34"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +5.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.45"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.69%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.20"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.26%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.960.61"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +5.80%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.54%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.03"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.49%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.99%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.553.71"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.32%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "434.75"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +6.06%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +4.95%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.59"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.51"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.72"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +10.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "38.54"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +5.27%  "

$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.90"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.57%  "

$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.08"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "709.95"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.53%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.39"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.08%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.86%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "42.24"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.80%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0844"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +27.68%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.55"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.93%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.87%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.37"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.35%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0476"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.88%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.08"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +7.31%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.70%  "

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.340"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.08%  "

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.24"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +8.16%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +7.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.54"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.23%  "

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "LidoDAOToken"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.46"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.67%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "148.98"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.21"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.89"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.33%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "25.34"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.75%  "
